$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUP_TRADE")

# Rename the "H2GC" transport blocks to "H2" (source cells that formulas below them reference)
$ws.Range("O4").Value = "TB_H2_DKISLBH_DKE_01"
$ws.Range("O9").Value = "TB_H2_DKISLBH_DKE_02"
$ws.Range("O14").Value = "TB_H2_DKISL1_DKW_01"
$ws.Range("O19").Value = "TB_H2_DKISL1_DKW_02"
$ws.Range("O24").Value = "TB_H2_DKISL2_DKW_01"
$ws.Range("O29").Value = "TB_H2_DKISL2_DKW_02"
$ws.Range("O34").Value = "TB_H2_DKISL3_DKW_01"
$ws.Range("O39").Value = "TB_H2_DKISL3_DKW_02"

# Update the active selection/view as in the author's saved state
$ws.Activate()
$ws.Range("O40").Select()
